$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the "A 5030-2019" record (currently row 9) down so it becomes
#        the last entry of the 2018/2019 block (final row 39), shifting the
#        rows in between (old rows 10-39) up by one. This mirrors an
#        Excel "copy row, insert copied row before the target, delete the
#        original row" sequence. ---
$ws.Rows.Item(9).Copy()
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(9).Delete()

# Re-apply the explicit row height on the row that now holds the moved
# record (row 39) - Insert() does not carry the custom row-height flag.
$ws.Rows.Item(39).RowHeight = 15

# --- 2. Bump the "Förändrad" (last-changed) date for every existing data
#        row from 45175 to 45177. ---
$ws.Range("C2:C301").Value = 45177

# --- 3. Append the new record "A 41515-2023" as row 302. ---
$ws.Range("A302").Value = "A 41515-2023"
$ws.Range("B302").Value = 45175
$ws.Range("B302").NumberFormat = "YYYY-MM-DD"
$ws.Range("C302").Value = 45177
$ws.Range("C302").NumberFormat = "YYYY-MM-DD"
$ws.Range("D302").Value = "ÖREBRO LÄN"
$ws.Range("E302").Value = "HALLSBERG"
$ws.Range("G302").Value = 2.4
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = 0
$ws.Range("J302").Value = 0
$ws.Range("K302").Value = 0
$ws.Range("L302").Value = 0
$ws.Range("M302").Value = 0
$ws.Range("N302").Value = 0
$ws.Range("O302").Value = 0
$ws.Range("P302").Value = 0
$ws.Range("Q302").Value = 0
$ws.Range("R302").WrapText = $true

# Row 301 is no longer the last row, so it regains the explicit row-height
# flag that Excel stamps on normal data rows (row 302, brand new, stays
# without it, matching a freshly appended row).
$ws.Rows.Item(301).RowHeight = 15
